$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13 (pushes "Programa resumido:" block and
# everything below it down by one row), then restore B/C number formatting on
# the new row by copying it from the row just below (which carries the
# correct column B/C styles) before clearing the inherited column-A cell.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 10 "Objetivos:" - new description text
$ws.Range("B10").Value = "Apresentar os princípios da automação da produção, características, aplicações e capacidades"
$ws.Range("C10").Value = "Apresentar os princípios da automação da produção, características, aplicações e capacidades"

# Row 13 (new) "Docentes responsáveis:" value, moved here from what is now row 19
$ws.Range("B13").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C13").Value = "8767640 - Eduardo Ferro dos Santos"

# Row 14 "Programa resumido:" - new short-syllabus text
$ws.Range("B14").Value = "Controle e automação; Robótica; Domótica; Sistemas Supervisórios, Pneumática, Hidráulica, CLP"
$ws.Range("C14").Value = "Controle e automação; Robótica; Domótica; Sistemas Supervisórios, Pneumática, Hidráulica, CLP"

# Row 16 "Programa:" - new full syllabus text
$ws.Range("B16").Value = "Introdução aos princípios de controle e automação; Fundamentos da Robótica; Fundamentos da Domótica;  Introdução a Sistemas Supervisórios, Princípios da Automação Pneumática, Hidráulica, Introdução aos Controladores Lógicos Programáveis."
$ws.Range("C16").Value = "Introdução aos princípios de controle e automação; Fundamentos da Robótica; Fundamentos da Domótica;  Introdução a Sistemas Supervisórios, Princípios da Automação Pneumática, Hidráulica, Introdução aos Controladores Lógicos Programáveis."

# Row 19 "Método:" - method text
$ws.Range("B19").Value = "Aulas expositivas e práticas."
$ws.Range("C19").Value = "Aulas expositivas e práticas."

# Row 20 "Critério:" - grading criteria text
$ws.Range("B20").Value = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"
$ws.Range("C20").Value = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"

# Row 21 "Norma de recuperação:" - recovery norm text
$ws.Range("B21").Value = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
$ws.Range("C21").Value = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."

# Row 22 "Bibliografia:" - new bibliography text
$ws.Range("B22").Value = "Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) NISE, N. S., “Engenharia de Sistemas de Controle”, 3ª ed., LTC, 2002. OGATA, K., “Engenharia de Controle Moderno”, 4ª ed., Prentice-Hall do Brasil, 2003. Tutoriais disponibilizados pelo professor BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U. B.. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p. CAPELLI, A. Automação Industrial: controle de movimento e processos contínuos. São Paulo: Érica, 2006. SILVEIRA, P. R. da; SANTOS, W. E. Automação e controle discreto. 3. ed. São Paulo: Érica, 1998. MORAES, C. C.; CATRUCCI, P. Engenharia de automação industrial. 2. ed. Rio de Janeiro: LTC, 2007. GIORGINI, M. Automação aplicada: descrição e implementação de sistemas sequencias com PLC's. 5. ed. São Paulo: Érica, 2003."
$ws.Range("C22").Value = "Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) NISE, N. S., “Engenharia de Sistemas de Controle”, 3ª ed., LTC, 2002. OGATA, K., “Engenharia de Controle Moderno”, 4ª ed., Prentice-Hall do Brasil, 2003. Tutoriais disponibilizados pelo professor BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U. B.. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p. CAPELLI, A. Automação Industrial: controle de movimento e processos contínuos. São Paulo: Érica, 2006. SILVEIRA, P. R. da; SANTOS, W. E. Automação e controle discreto. 3. ed. São Paulo: Érica, 1998. MORAES, C. C.; CATRUCCI, P. Engenharia de automação industrial. 2. ed. Rio de Janeiro: LTC, 2007. GIORGINI, M. Automação aplicada: descrição e implementação de sistemas sequencias com PLC's. 5. ed. São Paulo: Érica, 2003."
